$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "article 83 is live": the new blog post bumps the rotating blog slots in
# row 7 up by one. ser:80 drops off entirely, ser:81 moves from E7 into I7,
# ser:82 moves from C7 into E7, and the brand-new ser:83 takes over C7.
$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 81"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 82"
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 83"

# Leave the selection on I7, matching where the edit was made.
$ws.Range("I7").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 2
